$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{row=49; A='2021-04-05'; B='12:33:11'; C='2.233'; D='4598.6'; E='5254.4'; F='50.2'; G='51.6'; H='InService'},
    @{row=50; A='2021-04-05'; B='12:38:29'; C='2.5333'; D='4491.4'; E='6011.6'; F='51.2'; G='51.2'; H='InService'},
    @{row=51; A='2021-04-05'; B='12:38:57'; C='2.5333'; D='4491.4'; E='6011.6'; F='51.2'; G='51.2'; H='InService'},
    @{row=52; A='2021-04-05'; B='12:41:16'; C='2.5333'; D='4196.2'; E='5886.8'; F='43.0'; G='43.4'; H='InService'},
    @{row=53; A='2021-04-05'; B='12:41:32'; C='2.3333'; D='3360.0'; E='3000.0'; H='InService'},
    @{row=54; A='2021-04-05'; B='12:41:49'; C='2.3333'; D='3360.0'; E='3000.0'; H='InService'},
    @{row=55; A='2021-04-05'; B='15:43:04'; H='InService'},
    @{row=56; A='2021-04-05'; B='15:44:15'; H='InService'},
    @{row=57; A='2021-04-05'; B='16:28:34'; H='InService'},
    @{row=58; A='2021-04-05'; B='16:33:13'; C='4.8333'; D='2990.0'; E='2964.0'; H='InService'},
    @{row=59; A='2021-04-05'; B='16:33:38'; C='4.8333'; D='2990.0'; E='2964.0'; H='InService'},
    @{row=60; A='2021-04-05'; B='16:37:00'; C='0.5'; D='6634.0'; E='7839.0'; H='InService'},
    @{row=61; A='2021-04-05'; B='16:40:07'; C='2.3333'; D='3392.0'; E='2952.0'; H='InService'},
    @{row=62; A='2021-04-05'; B='16:41:04'; C='2.2951'; D='3360.0'; E='3100.0'; H='InService'},
    @{row=63; A='2021-04-05'; B='16:42:25'; C='2.3729'; D='3426.0'; E='3100.0'; H='InService'},
    @{row=64; A='2021-04-05'; B='16:43:24'; C='20.8333'; D='5399.0'; E='3644.0'; H='InService'},
    @{row=65; A='2021-04-05'; B='16:43:49'; C='20.8333'; D='5399.0'; E='3644.0'; H='InService'},
    @{row=66; A='2021-04-05'; B='16:44:56'; C='25.1667'; D='5764.0'; E='5536.0'; H='InService'},
    @{row=67; A='2021-04-05'; B='16:45:56'; C='27.3333'; D='6965.0'; E='5536.0'; H='InService'},
    @{row=68; A='2021-04-05'; B='16:47:05'; C='29.5'; D='7989.0'; E='6688.0'; H='InService'},
    @{row=69; A='2021-04-05'; B='16:48:18'; C='6.0656'; D='4828.0'; E='5196.0'; H='InService'},
    @{row=70; A='2021-04-05'; B='16:49:24'; C='0.5'; D='2536.0'; E='2444.0'; H='InService'},
    @{row=71; A='2021-04-05'; B='16:50:24'; C='0.5'; D='2674.0'; E='2312.0'; H='InService'},
    @{row=72; A='2021-04-05'; B='16:51:24'; C='0.5'; D='2682.0'; E='2548.0'; H='InService'},
    @{row=73; A='2021-04-05'; B='16:52:24'; C='0.5'; D='4950.0'; E='4937.0'; H='InService'},
    @{row=74; A='2021-04-05'; B='16:53:24'; C='2.3333'; D='3450.0'; E='3176.0'; H='InService'},
    @{row=75; A='2021-04-05'; B='16:57:55'; C='2.7869'; D='3930.0'; E='3948.0'; H='InService'},
    @{row=76; A='2021-04-05'; B='16:59:02'; C='2.2951'; D='5751.0'; E='5739.0'; H='InService'},
    @{row=77; A='2021-04-05'; B='17:00:02'; C='21.3115'; D='5695.0'; E='3984.0'; H='InService'},
    @{row=78; A='2021-04-05'; B='17:01:02'; C='31.5517'; D='8878.0'; E='7791.0'; H='InService'},
    @{row=79; A='2021-04-05'; B='17:02:02'; C='29.0323'; D='7023.0'; E='5476.0'; H='InService'},
    @{row=80; A='2021-04-05'; B='17:03:02'; C='12.931'; D='5137.0'; E='4804.0'; H='InService'},
    @{row=81; A='2021-04-05'; B='17:04:02'; C='20.8197'; D='11172.0'; E='10636.0'; H='InService'},
    @{row=82; A='2021-04-05'; B='17:05:02'; C='21.3333'; D='6219.0'; E='5836.0'; H='InService'}
)

foreach ($item in $rows) {
    $r = $item.row
    foreach ($col in @("A","B","C","D","E","F","G","H")) {
        if ($item.ContainsKey($col)) {
            $cellRef = "$col$r"
            $ws.Range($cellRef).NumberFormat = "@"
            $ws.Range($cellRef).Value = $item[$col]
            $ws.Range($cellRef).Style = "Normal"
        }
    }
}
